$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 578
$ws1.Range("F5").Value = 1815
$ws1.Range("F6").Value = 279
$ws1.Range("F8").Value = 162
$ws1.Range("F9").Value = 2325
$ws1.Range("F10").Value = 123
$ws1.Range("F13").Value = 1414
$ws1.Range("F14").Value = 501
$ws1.Range("F17").Value = 219
$ws1.Range("F18").Value = 15
$ws1.Range("F24").Value = 80
$ws1.Range("F25").Value = 31
$ws1.Range("F26").Value = 1447
$ws1.Range("F28").Value = 367
$ws1.Range("F29").Value = 198
$ws1.Range("F30").Value = 182
$ws1.Range("F31").Value = 285
$ws1.Range("F32").Value = 360

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 578
$ws4.Range("F5").Value = 1815
$ws4.Range("F7").Value = 279
$ws4.Range("F9").Value = 162
$ws4.Range("F10").Value = 2325
$ws4.Range("F11").Value = 123
$ws4.Range("F14").Value = 1414
$ws4.Range("F15").Value = 501
$ws4.Range("F18").Value = 219
$ws4.Range("F19").Value = 15
$ws4.Range("F25").Value = 80
$ws4.Range("F26").Value = 31
$ws4.Range("F27").Value = 1447
$ws4.Range("F29").Value = 367
$ws4.Range("F30").Value = 198
$ws4.Range("F31").Value = 182
$ws4.Range("F32").Value = 285
$ws4.Range("F33").Value = 360
